# Updated cryptos list on Tue Sep 26 14:00:24 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns for the existing coins
# and inserts a new coin (BabyDogeCoin) at the bottom of the ranking, which
# pushes RenderToken / Algorand / Cronos down a row and drops the former
# last row (USDD) off the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (URLs, percentages, coin names, and Price values that
# Excel can't mistake for a number because they contain two '.' separators
# or non-numeric characters).
$textUpdates = @{
    "D2"  = "26.206.30"
    "E2"  = "  +0.12%  "
    "D3"  = "1.585.83"
    "E3"  = "  +0.59%  "
    "E4"  = "  +0.11%  "
    "E5"  = "  +1.49%  "
    "E6"  = "  +0.78%  "
    "E7"  = "  +0.02%  "
    "E8"  = "  +0.53%  "
    "E9"  = "  -0.39%  "
    "E10" = "  -1.45%  "
    "E11" = "  +0.56%  "
    "D12" = "1.810.98"
    "E12" = "  +0.73%  "
    "D13" = "1.581.03"
    "E13" = "  +0.00%  "
    "E14" = "  -1.21%  "
    "E15" = "  +0.72%  "
    "E16" = "  -0.43%  "
    "D17" = "26.238.36"
    "E17" = "  +0.30%  "
    "E18" = "  -0.19%  "
    "E19" = "  +1.37%  "
    "E20" = "  +2.10%  "
    "E21" = "  -0.06%  "
    "E22" = "  -0.12%  "
    "E23" = "  -0.25%  "
    "E24" = "  +1.55%  "
    "E26" = "  -0.01%  "
    "E27" = "  +0.13%  "
    "E28" = "  -0.75%  "
    "E29" = "  -0.27%  "
    "E30" = "  -2.10%  "
    "E31" = "  +1.30%  "
    "E32" = "  -0.48%  "
    "D33" = "1.338.58"
    "E33" = "  +4.74%  "
    "E34" = "  -2.41%  "
    "E35" = "  -0.09%  "
    "E36" = "  -1.18%  "
    "E37" = "  -5.34%  "
    "E38" = "  +0.53%  "
    "E39" = "  +1.37%  "
    "E40" = "  +3.33%  "
    "E41" = "  -0.02%  "
    "E42" = "  -14.33%  "
    "E43" = "  +0.52%  "
    "E44" = "  +0.51%  "
    "D45" = "1.723.52"
    "E45" = "  +0.74%  "
    "E46" = "  -2.04%  "
    "E47" = "  -3.48%  "

    "B48" = "BabyDogeCoin"
    "C48" = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
    "D48" = "0.0₆0103"
    "E48" = "  -0.48%  "

    "B49" = "RenderToken"
    "C49" = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
    "E49" = "  -1.56%  "

    "B50" = "Algorand"
    "C50" = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
    "E50" = "  -2.74%  "

    "B51" = "Cronos"
    "C51" = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
    "E51" = "  -0.88%  "
}

foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

# Price values that look like plain numbers to Excel (single '.' decimal
# point) need to be forced to text, the same way the source data keeps them
# as strings (e.g. so "19.29" isn't silently turned into the number 19.29).
# Temporarily switch the cell to a text format, assign the value, then
# restore the row's normal (unstyled) look from an untouched neighbour cell
# so no stray number format is left behind.
$numericLookingPrices = @{
    "D5"  = "211.98"
    "D10" = "19.29"
    "D16" = "64.15"
    "D19" = "7.36"
    "D22" = "4.25"
    "D23" = "2.16"
    "D24" = "8.98"
    "D25" = "143.54"
    "D29" = "15.17"
    "D30" = "0.0495"
    "D31" = "1.16"
    "D34" = "2.94"
    "D37" = "0.580"
    "D39" = "0.822"
    "D42" = "0.952"
    "D46" = "61.11"
    "D47" = "85.69"
    "D49" = "1.48"
    "D50" = "0.0978"
    "D51" = "0.0501"
}

foreach ($addr in $numericLookingPrices.Keys) {
    $row = $addr.Substring(1)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $numericLookingPrices[$addr]
    $ws.Range($addr).Style = $ws.Range("B" + $row).Style
}
